$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 2584
$ws1.Range("F3").Value = 561
$ws1.Range("F4").Value = 461
$ws1.Range("F7").Value = 461
$ws1.Range("F8").Value = 1199
$ws1.Range("F9").Value = 540
$ws1.Range("F10").Value = 299
$ws1.Range("F12").Value = 347
$ws1.Range("F13").Value = 5582
$ws1.Range("F14").Value = 60
$ws1.Range("F15").Value = 1704
$ws1.Range("F16").Value = 4031
$ws1.Range("F17").Value = 409
$ws1.Range("F20").Value = 4632
$ws1.Range("F21").Value = 6050
$ws1.Range("F23").Value = 1033
$ws1.Range("F24").Value = 667
$ws1.Range("F25").Value = 3708
$ws1.Range("F26").Value = 483
$ws1.Range("F28").Value = 183
$ws1.Range("F29").Value = 121
$ws1.Range("F30").Value = 973
$ws1.Range("F31").Value = 1374
$ws1.Range("F32").Value = 455
$ws1.Range("F33").Value = 521
$ws1.Range("F34").Value = 1565
$ws1.Range("F35").Value = 194
$ws1.Range("F36").Value = 1664
$ws1.Range("F37").Value = 167
$ws1.Range("F39").Value = 1100
$ws1.Range("F40").Value = 31
$ws1.Range("F42").Value = 611
$ws1.Range("F43").Value = 87
$ws1.Range("F44").Value = 216
$ws1.Range("F45").Value = 3312
$ws1.Range("F47").Value = 272
$ws1.Range("F48").Value = 402
$ws1.Range("F49").Value = 3865

# Sheet: 演出
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 1179

# Sheet: 本地生活
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 3730

# Sheet: 全部类型
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 3730
$ws4.Range("F3").Value = 2584
$ws4.Range("F4").Value = 561
$ws4.Range("F5").Value = 461
$ws4.Range("F7").Value = 1179
$ws4.Range("F10").Value = 461
$ws4.Range("F11").Value = 1199
$ws4.Range("F12").Value = 540
$ws4.Range("F13").Value = 299
$ws4.Range("F15").Value = 347
$ws4.Range("F17").Value = 1704
$ws4.Range("F18").Value = 4632
$ws4.Range("F19").Value = 6050
$ws4.Range("F21").Value = 1033
$ws4.Range("F22").Value = 667
$ws4.Range("F23").Value = 3708
$ws4.Range("F24").Value = 484
$ws4.Range("F26").Value = 183
$ws4.Range("F27").Value = 121
$ws4.Range("F28").Value = 973
$ws4.Range("F29").Value = 1374
$ws4.Range("F30").Value = 455
$ws4.Range("F31").Value = 521
$ws4.Range("F33").Value = 1565
$ws4.Range("F34").Value = 194
$ws4.Range("F35").Value = 1664
$ws4.Range("F37").Value = 1100
$ws4.Range("F39").Value = 611
$ws4.Range("F41").Value = 87
$ws4.Range("F43").Value = 3312
$ws4.Range("F46").Value = 272
$ws4.Range("F47").Value = 402
$ws4.Range("F49").Value = 3865

"done"